$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# For each changed Price/Volume cell: force Text format, write the new
# literal value (so Excel will not re-interpret numeric-looking strings
# like '0.9987' or '1.0000' as real numbers), then clear the temporary
# Text number-format so the cell's style index matches the original file
# (plain/default style, just like the un-touched cells around it).
# NumberFormat/ClearFormats are applied per single cell (not as a multi-
# area union range) to avoid a COM quirk that only formats the first area.
function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

Set-TextValue 'D2' '25.761.62'
Set-TextValue 'E2' '  +3.17%  '
Set-TextValue 'D3' '1.676.46'
Set-TextValue 'E3' '  +2.20%  '
Set-TextValue 'D4' '0.9987'
Set-TextValue 'E4' '  -0.04%  '
Set-TextValue 'D5' '237.42'
Set-TextValue 'E5' '  +2.10%  '
Set-TextValue 'D6' '1.0000'
Set-TextValue 'E6' '  -0.02%  '
Set-TextValue 'D7' '0.4628'
Set-TextValue 'E7' '  -2.55%  '
Set-TextValue 'E8' '  +0.50%  '
Set-TextValue 'D9' '0.06142'
Set-TextValue 'E9' '  +0.96%  '
Set-TextValue 'D10' '1.672.83'
Set-TextValue 'E10' '  +1.99%  '
Set-TextValue 'D11' '0.06994'
Set-TextValue 'E11' '  -0.43%  '
Set-TextValue 'E12' '  +2.65%  '
Set-TextValue 'E13' '  +1.03%  '
Set-TextValue 'D14' '0.5766'
Set-TextValue 'E14' '  -1.36%  '
Set-TextValue 'E15' '  +2.33%  '
Set-TextValue 'D16' '0.9998'
Set-TextValue 'E16' '  -0.02%  '
Set-TextValue 'D17' '0.9998'
Set-TextValue 'E17' '  -0.01%  '
Set-TextValue 'D18' '25.759.11'
Set-TextValue 'E18' '  +3.17%  '
Set-TextValue 'E19' '  +1.71%  '
Set-TextValue 'D20' '11.43'
Set-TextValue 'E20' '  +1.71%  '
Set-TextValue 'D21' '1.886.14'
Set-TextValue 'E21' '  +1.54%  '
Set-TextValue 'D22' '4.466'
Set-TextValue 'E22' '  +3.64%  '
Set-TextValue 'D23' '8.667'
Set-TextValue 'E23' '  +1.52%  '
Set-TextValue 'D24' '5.227'
Set-TextValue 'E24' '  +0.00%  '
Set-TextValue 'D25' '134.13'
Set-TextValue 'E25' '  +0.53%  '
Set-TextValue 'D26' '14.95'
Set-TextValue 'E26' '  +0.53%  '
Set-TextValue 'E27' '  +0.80%  '
Set-TextValue 'D28' '1.717'
Set-TextValue 'E28' '  +5.07%  '
Set-TextValue 'D29' '104.45'
Set-TextValue 'E29' '  +0.25%  '
Set-TextValue 'D30' '3.947'
Set-TextValue 'E30' '  +1.24%  '
Set-TextValue 'D31' '0.07686'
Set-TextValue 'E31' '  +1.66%  '
Set-TextValue 'D32' '3.616'
Set-TextValue 'E32' '  +1.51%  '
Set-TextValue 'D33' '0.04338'
Set-TextValue 'E33' '  +1.53%  '
Set-TextValue 'D34' '2.599'
Set-TextValue 'E34' '  +1.10%  '
Set-TextValue 'D35' '0.6074'
Set-TextValue 'E35' '  +2.36%  '
Set-TextValue 'D36' '0.9501'
Set-TextValue 'E36' '  +2.38%  '
Set-TextValue 'D37' '0.9343'
Set-TextValue 'E37' '  +4.66%  '
Set-TextValue 'D38' '109.03'
Set-TextValue 'E38' '  +10.53%  '
Set-TextValue 'D39' '2.445'
Set-TextValue 'E39' '  -5.18%  '
Set-TextValue 'D40' '0.9990'
Set-TextValue 'E40' '  -0.05%  '
Set-TextValue 'D41' '1.858'
Set-TextValue 'E41' '  +5.44%  '
Set-TextValue 'D42' '0.01450'
Set-TextValue 'E42' '  -3.06%  '
Set-TextValue 'D43' '5.046'
Set-TextValue 'E43' '  +8.29%  '
Set-TextValue 'D44' '0.3723'
Set-TextValue 'E44' '  +0.46%  '
Set-TextValue 'D45' '0.1119'
Set-TextValue 'E45' '  +1.53%  '
Set-TextValue 'D46' '0.05301'
Set-TextValue 'E46' '  +2.00%  '
Set-TextValue 'D47' '31.32'
Set-TextValue 'E47' '  +9.31%  '
Set-TextValue 'D48' '6.146'
Set-TextValue 'E48' '  +0.82%  '
Set-TextValue 'D49' '7.616'
Set-TextValue 'E49' '  +7.01%  '
Set-TextValue 'D50' '1.212'
Set-TextValue 'E50' '  +2.56%  '
Set-TextValue 'E51' '  +0.00%  '
